$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 6-9 (data for MuSCs/Resolving-Mac -> ECs/Resolving-Mac pairs
# collapsed away with the new TPM numbers), shifting nothing else below them.
$ws.Rows("6:9").Delete()

# Row 2: ECs -> ECs/Icam1/Itgam -> Resolving-Mac
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam1"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 30.87085333333333
$ws.Range("H2").Value = 92.61256
$ws.Range("I2").Value = 0.2985789950947061
$ws.Range("J2").Value = 0.2985789950947061
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 35.68243999999999
$ws.Range("N2").Value = 107.04732
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1101.547371815466
$ws.Range("R2").Value = 9913.926346339198
$ws.Range("S2").Value = 0.2985789950947061
$ws.Range("T2").Value = 0.2985789950947061

# Row 3: ECs -> FAPs/Icam1/Itgam -> Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Icam1"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.793597
$ws.Range("H3").Value = 101.380791
$ws.Range("I3").Value = 0.3268474027571036
$ws.Range("J3").Value = 0.3268474027571037
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 35.68243999999999
$ws.Range("N3").Value = 107.04732
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1205.83799733668
$ws.Range("R3").Value = 10852.54197603012
$ws.Range("S3").Value = 0.3268474027571036
$ws.Range("T3").Value = 0.3268474027571037

# Row 4: FAPs -> MuSCs/Icam1/Itgam -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Icam1"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.981185666666667
$ws.Range("H4").Value = 8.943557
$ws.Range("I4").Value = 0.02883365130639111
$ws.Range("J4").Value = 0.02883365130639111
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 35.68243999999999
$ws.Range("N4").Value = 107.04732
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 106.3759786796933
$ws.Range("R4").Value = 957.3838081172399
$ws.Range("S4").Value = 0.02883365130639111
$ws.Range("T4").Value = 0.02883365130639111

# Row 5: FAPs -> Resolving-Mac/Icam1/Itgam -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Icam1"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 35.74694633333333
$ws.Range("H5").Value = 107.240839
$ws.Range("I5").Value = 0.3457399508417991
$ws.Range("J5").Value = 0.3457399508417991
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 35.68243999999999
$ws.Range("N5").Value = 107.04732
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1275.538267722386
$ws.Range("R5").Value = 11479.84440950148
$ws.Range("S5").Value = 0.3457399508417991
$ws.Range("T5").Value = 0.3457399508417991
